$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value2 = 14494481
$ws.Range("I106").Value2 = 17545166
$ws.Range("J106").Value2 = 3727.5
$ws.Range("K106").Value2 = 17545166
$ws.Range("L106").Value2 = 3727.5
$ws.Range("M106").Value2 = -17544535
$ws.Range("N106").Value2 = -4989.5
$ws.Range("H129").Value2 = 170506.84
$ws.Range("J129").Value2 = 176472
$ws.Range("L129").Value2 = 529416
$ws.Range("N129").Value2 = -539416
$ws.Range("H132").Value2 = 2004.2157
$ws.Range("I132").Value2 = 2091.9148
$ws.Range("K132").Value2 = 6275.7444
$ws.Range("M132").Value2 = -3745.7444
$ws.Range("H135").Value2 = 13893585
$ws.Range("I135").Value2 = 453
$ws.Range("J135").Value2 = 100031000
$ws.Range("K135").Value2 = 4077
$ws.Range("L135").Value2 = 900279000
$ws.Range("M135").Value2 = -1542
$ws.Range("N135").Value2 = -900284070
$ws.Range("H138").Value2 = 1795.5784
$ws.Range("I138").Value2 = 759
$ws.Range("J138").Value2 = 2242.3794
$ws.Range("K138").Value2 = 2277
$ws.Range("L138").Value2 = 6727.138199999999
$ws.Range("M138").Value2 = 2863
$ws.Range("N138").Value2 = -17007.1382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 17305.3
$ws.Range("I32").Value2 = 19332.766
$ws.Range("J32").Value2 = 2887.7778
$ws.Range("K32").Value2 = 19332.766
$ws.Range("L32").Value2 = 2887.7778
$ws.Range("M32").Value2 = -19045.766
$ws.Range("N32").Value2 = -3461.7778
$ws.Range("H61").Value2 = 1777.8975
$ws.Range("I61").Value2 = 1549.3103
$ws.Range("K61").Value2 = 1549.3103
$ws.Range("M61").Value2 = -1337.3103
$ws.Range("H97").Value2 = 1192.7307
$ws.Range("I97").Value2 = 1356.0555
$ws.Range("K97").Value2 = 1356.0555
$ws.Range("M97").Value2 = -860.0554999999999
$ws.Range("H122").Value2 = 2554.8823
$ws.Range("I122").Value2 = 1640.3334
$ws.Range("K122").Value2 = 4921.0002
$ws.Range("M122").Value2 = -2471.0002
$ws.Range("H132").Value2 = 13557.675
$ws.Range("I132").Value2 = 1853.1818
$ws.Range("J132").Value2 = 52182.5
$ws.Range("K132").Value2 = 5559.5454
$ws.Range("L132").Value2 = 156547.5
$ws.Range("M132").Value2 = -3029.5454
$ws.Range("N132").Value2 = -161607.5
$ws.Range("H136").Value2 = 1777.8975
$ws.Range("I136").Value2 = 1549.3103
$ws.Range("K136").Value2 = 4647.9309
$ws.Range("M136").Value2 = -2097.9309

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 69300
$ws.Range("I134").Value2 = 91800.086
$ws.Range("K134").Value2 = 275400.258
$ws.Range("M134").Value2 = -272865.258

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 10235.475
$ws.Range("I31").Value2 = 21282.125
$ws.Range("J31").Value2 = 2871.0417
$ws.Range("K31").Value2 = 21282.125
$ws.Range("L31").Value2 = 2871.0417
$ws.Range("M31").Value2 = -20987.125
$ws.Range("N31").Value2 = -3461.0417
$ws.Range("H34").Value2 = 10235.475
$ws.Range("I34").Value2 = 21282.125
$ws.Range("J34").Value2 = 2871.0417
$ws.Range("K34").Value2 = 21282.125
$ws.Range("L34").Value2 = 2871.0417
$ws.Range("M34").Value2 = -21080.125
$ws.Range("N34").Value2 = -3275.0417
$ws.Range("H56").Value2 = 18666.666
$ws.Range("J56").Value2 = 18666.666
$ws.Range("L56").Value2 = 18666.666
$ws.Range("N56").Value2 = -20356.666
$ws.Range("H58").Value2 = 12985.342
$ws.Range("I58").Value2 = 743.4286
$ws.Range("J58").Value2 = 39352.54
$ws.Range("K58").Value2 = 743.4286
$ws.Range("L58").Value2 = 39352.54
$ws.Range("M58").Value2 = -540.4286
$ws.Range("N58").Value2 = -39758.54
$ws.Range("H86").Value2 = 11919749
$ws.Range("I86").Value2 = 7806.125
$ws.Range("K86").Value2 = 7806.125
$ws.Range("M86").Value2 = -6683.125
$ws.Range("H89").Value2 = 11919749
$ws.Range("I89").Value2 = 7806.125
$ws.Range("K89").Value2 = 39030.625
$ws.Range("M89").Value2 = -33414.625
$ws.Range("H99").Value2 = 16132700
$ws.Range("I99").Value2 = 3668.75
$ws.Range("K99").Value2 = 3668.75
$ws.Range("M99").Value2 = -2170.75
$ws.Range("H122").Value2 = 1039
$ws.Range("I122").Value2 = 992.2273
$ws.Range("J122").Value2 = 1081.875
$ws.Range("K122").Value2 = 2976.6819
$ws.Range("L122").Value2 = 3245.625
$ws.Range("M122").Value2 = -526.6819
$ws.Range("N122").Value2 = -8145.625
$ws.Range("H126").Value2 = 16132700
$ws.Range("I126").Value2 = 3668.75
$ws.Range("K126").Value2 = 11006.25
$ws.Range("M126").Value2 = -8536.25
$ws.Range("H134").Value2 = 1033.5333
$ws.Range("I134").Value2 = 958.5833
$ws.Range("J134").Value2 = 1333.3334
$ws.Range("K134").Value2 = 2875.7499
$ws.Range("L134").Value2 = 4000.0002
$ws.Range("M134").Value2 = -340.7498999999998
$ws.Range("N134").Value2 = -9070.0002
$ws.Range("H136").Value2 = 12985.342
$ws.Range("I136").Value2 = 743.4286
$ws.Range("J136").Value2 = 39352.54
$ws.Range("K136").Value2 = 2230.2858
$ws.Range("L136").Value2 = 118057.62
$ws.Range("M136").Value2 = 319.7142000000003
$ws.Range("N136").Value2 = -123157.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value2 = 973.5454999999999
$ws.Range("I18").Value2 = 141.75
$ws.Range("J18").Value2 = 3191.6667
$ws.Range("K18").Value2 = 425.25
$ws.Range("L18").Value2 = 9575.000100000001
$ws.Range("M18").Value2 = -256.25
$ws.Range("N18").Value2 = -9913.000100000001
$ws.Range("H98").Value2 = 824.6667
$ws.Range("I98").Value2 = 1899
$ws.Range("J98").Value2 = 287.5
$ws.Range("K98").Value2 = 5697
$ws.Range("L98").Value2 = 862.5
$ws.Range("M98").Value2 = -4199
$ws.Range("N98").Value2 = -3858.5
$ws.Range("H131").Value2 = 756.47
$ws.Range("I131").Value2 = 223.25
$ws.Range("J131").Value2 = 802.837
$ws.Range("K131").Value2 = 669.75
$ws.Range("L131").Value2 = 2408.511
$ws.Range("M131").Value2 = 4370.25
$ws.Range("N131").Value2 = -12488.511

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 41668116
$ws.Range("I102").Value2 = 45455900
$ws.Range("J102").Value2 = 2500
$ws.Range("K102").Value2 = 45455900
$ws.Range("L102").Value2 = 2500
$ws.Range("M102").Value2 = -45454278
$ws.Range("N102").Value2 = -5744
$ws.Range("H107").Value2 = 5494703
$ws.Range("I107").Value2 = 198.09091
$ws.Range("J107").Value2 = 25641220
$ws.Range("K107").Value2 = 198.09091
$ws.Range("L107").Value2 = 25641220
$ws.Range("M107").Value2 = 1721.90909
$ws.Range("N107").Value2 = -25645060
$ws.Range("H113").Value2 = 2810.0286
$ws.Range("I113").Value2 = 2648.4092
$ws.Range("J113").Value2 = 3083.5386
$ws.Range("K113").Value2 = 2648.4092
$ws.Range("L113").Value2 = 3083.5386
$ws.Range("M113").Value2 = -478.4092000000001
$ws.Range("N113").Value2 = -7423.5386
$ws.Range("H122").Value2 = 45978664
$ws.Range("I122").Value2 = 18519844
$ws.Range("J122").Value2 = 90911280
$ws.Range("K122").Value2 = 55559532
$ws.Range("L122").Value2 = 272733840
$ws.Range("M122").Value2 = -55557082
$ws.Range("N122").Value2 = -272738740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 6175
$ws.Range("I7").Value2 = 4780
$ws.Range("J7").Value2 = 6809.091
$ws.Range("K7").Value2 = 4780
$ws.Range("L7").Value2 = 6809.091
$ws.Range("M7").Value2 = -4668
$ws.Range("N7").Value2 = -7033.091
$ws.Range("H40").Value2 = 6869
$ws.Range("I40").Value2 = 5199.8
$ws.Range("J40").Value2 = 8260
$ws.Range("K40").Value2 = 5199.8
$ws.Range("L40").Value2 = 8260
$ws.Range("M40").Value2 = -5063.8
$ws.Range("N40").Value2 = -8532
$ws.Range("H95").Value2 = 40000
$ws.Range("J95").Value2 = 40000
$ws.Range("L95").Value2 = 40000
$ws.Range("N95").Value2 = -45492
$ws.Range("H97").Value2 = 19672
$ws.Range("J97").Value2 = 19672
$ws.Range("L97").Value2 = 19672
$ws.Range("N97").Value2 = -21654
$ws.Range("H126").Value2 = 6175
$ws.Range("I126").Value2 = 4780
$ws.Range("J126").Value2 = 6809.091
$ws.Range("K126").Value2 = 14340
$ws.Range("L126").Value2 = 20427.273
$ws.Range("M126").Value2 = -11870
$ws.Range("N126").Value2 = -25367.273
$ws.Range("H132").Value2 = 1903.7916
$ws.Range("I132").Value2 = 1420.1333
$ws.Range("J132").Value2 = 2709.889
$ws.Range("K132").Value2 = 4260.3999
$ws.Range("L132").Value2 = 8129.667
$ws.Range("M132").Value2 = -1730.3999
$ws.Range("N132").Value2 = -13189.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value2 = 18000
$ws.Range("I61").Value2 = 0
$ws.Range("K61").Value2 = 0
$ws.Range("M61").ClearContents()
$ws.Range("H107").Value2 = 3030986.5
$ws.Range("J107").Value2 = 7576257.5
$ws.Range("L107").Value2 = 22728772.5
$ws.Range("N107").Value2 = -22732612.5
$ws.Range("H122").Value2 = 2156.4285
$ws.Range("I122").Value2 = 2016
$ws.Range("K122").Value2 = 6048
$ws.Range("M122").Value2 = -3598

Write-Output "Applied 222 cell edits across 8 sheets"